$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "84ec514e-b81d-4773-8859-f85d7adf8d5d"
$ws.Range("B8").Value = "Create Citizenship by getting data from Excel"
$ws.Range("C8").Value = "FAILED"
$ws.Range("D8").Value = "2023-10-05T03:08:46.196461"
$ws.Range("E8").Value = "2023-10-05T03:09:01.847120500"
$ws.Range("F8").Value = "PT15.6506595S"
